$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

# Row 45
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = "2026-01-28"
$ws.Range("A45").Style = "Normal"
$ws.Range("B45").Value = "17:23:06"
$ws.Range("C45").Value = "17:00"
$ws.Range("D45").Value = "Bedroom"
$ws.Range("E45").Value = "In Bed | HR=0 | BR=0"
$ws.Range("F45").Value = "Occupied"

# Row 46
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = "2026-01-28"
$ws.Range("A46").Style = "Normal"
$ws.Range("B46").Value = "17:23:07"
$ws.Range("C46").Value = "17:00"
$ws.Range("D46").Value = "Bedroom"
$ws.Range("E46").Value = "In Bed | HR=104 | BR=56"
$ws.Range("F46").Value = "Occupied"

# Row 47
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = "2026-01-28"
$ws.Range("A47").Style = "Normal"
$ws.Range("B47").Value = "17:23:08"
$ws.Range("C47").Value = "17:00"
$ws.Range("D47").Value = "Bedroom"
$ws.Range("E47").Value = "In Bed | HR=61 | BR=13"
$ws.Range("F47").Value = "Occupied"

# Row 48
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = "2026-01-28"
$ws.Range("A48").Style = "Normal"
$ws.Range("B48").Value = "17:23:09"
$ws.Range("C48").Value = "17:00"
$ws.Range("D48").Value = "Bedroom"
$ws.Range("E48").Value = "In Bed | HR=54 | BR=6"
$ws.Range("F48").Value = "Occupied"

# Row 49
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = "2026-01-28"
$ws.Range("A49").Style = "Normal"
$ws.Range("B49").Value = "17:23:11"
$ws.Range("C49").Value = "17:00"
$ws.Range("D49").Value = "Bedroom"
$ws.Range("E49").Value = "In Bed | HR=69 | BR=21"
$ws.Range("F49").Value = "Occupied"

# Row 50
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "2026-01-28"
$ws.Range("A50").Style = "Normal"
$ws.Range("B50").Value = "17:23:12"
$ws.Range("C50").Value = "17:00"
$ws.Range("D50").Value = "Bedroom"
$ws.Range("E50").Value = "In Bed | HR=50 | BR=2"
$ws.Range("F50").Value = "Occupied"

# Row 51
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = "2026-01-28"
$ws.Range("A51").Style = "Normal"
$ws.Range("B51").Value = "17:23:16"
$ws.Range("C51").Value = "17:00"
$ws.Range("D51").Value = "Bedroom"
$ws.Range("E51").Value = "In Bed | HR=94 | BR=46"
$ws.Range("F51").Value = "Occupied"

# Row 52
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = "2026-01-28"
$ws.Range("A52").Style = "Normal"
$ws.Range("B52").Value = "17:23:17"
$ws.Range("C52").Value = "17:00"
$ws.Range("D52").Value = "Bedroom"
$ws.Range("E52").Value = "In Bed | HR=50 | BR=2"
$ws.Range("F52").Value = "Occupied"

# Row 53
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "2026-01-28"
$ws.Range("A53").Style = "Normal"
$ws.Range("B53").Value = "17:23:20"
$ws.Range("C53").Value = "17:00"
$ws.Range("D53").Value = "Bedroom"
$ws.Range("E53").Value = "In Bed | HR=49 | BR=1"
$ws.Range("F53").Value = "Occupied"

# Row 54
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value = "2026-01-28"
$ws.Range("A54").Style = "Normal"
$ws.Range("B54").Value = "17:23:35"
$ws.Range("C54").Value = "17:00"
$ws.Range("D54").Value = "Bedroom"
$ws.Range("E54").Value = "In Bed | HR=122 | BR=74"
$ws.Range("F54").Value = "Occupied"

# Row 55
$ws.Range("A55").NumberFormat = "@"
$ws.Range("A55").Value = "2026-01-28"
$ws.Range("A55").Style = "Normal"
$ws.Range("B55").Value = "17:23:36"
$ws.Range("C55").Value = "17:00"
$ws.Range("D55").Value = "Bedroom"
$ws.Range("E55").Value = "In Bed | HR=59 | BR=11"
$ws.Range("F55").Value = "Occupied"

# Row 56
$ws.Range("A56").NumberFormat = "@"
$ws.Range("A56").Value = "2026-01-28"
$ws.Range("A56").Style = "Normal"
$ws.Range("B56").Value = "17:23:37"
$ws.Range("C56").Value = "17:00"
$ws.Range("D56").Value = "Bedroom"
$ws.Range("E56").Value = "In Bed | HR=76 | BR=28"
$ws.Range("F56").Value = "Occupied"

# Row 57
$ws.Range("A57").NumberFormat = "@"
$ws.Range("A57").Value = "2026-01-28"
$ws.Range("A57").Style = "Normal"
$ws.Range("B57").Value = "17:23:38"
$ws.Range("C57").Value = "17:00"
$ws.Range("D57").Value = "Bedroom"
$ws.Range("E57").Value = "In Bed | HR=104 | BR=56"
$ws.Range("F57").Value = "Occupied"

# Row 58
$ws.Range("A58").NumberFormat = "@"
$ws.Range("A58").Value = "2026-01-28"
$ws.Range("A58").Style = "Normal"
$ws.Range("B58").Value = "17:23:39"
$ws.Range("C58").Value = "17:00"
$ws.Range("D58").Value = "Bedroom"
$ws.Range("E58").Value = "In Bed | HR=112 | BR=64"
$ws.Range("F58").Value = "Occupied"

# Row 59
$ws.Range("A59").NumberFormat = "@"
$ws.Range("A59").Value = "2026-01-28"
$ws.Range("A59").Style = "Normal"
$ws.Range("B59").Value = "17:23:40"
$ws.Range("C59").Value = "17:00"
$ws.Range("D59").Value = "Bedroom"
$ws.Range("E59").Value = "In Bed | HR=58 | BR=10"
$ws.Range("F59").Value = "Occupied"

# Row 60
$ws.Range("A60").NumberFormat = "@"
$ws.Range("A60").Value = "2026-01-28"
$ws.Range("A60").Style = "Normal"
$ws.Range("B60").Value = "17:23:41"
$ws.Range("C60").Value = "17:00"
$ws.Range("D60").Value = "Bedroom"
$ws.Range("E60").Value = "In Bed | HR=50 | BR=2"
$ws.Range("F60").Value = "Occupied"

# Row 61
$ws.Range("A61").NumberFormat = "@"
$ws.Range("A61").Value = "2026-01-28"
$ws.Range("A61").Style = "Normal"
$ws.Range("B61").Value = "17:23:42"
$ws.Range("C61").Value = "17:00"
$ws.Range("D61").Value = "Bedroom"
$ws.Range("E61").Value = "In Bed | HR=54 | BR=6"
$ws.Range("F61").Value = "Occupied"

# Row 62
$ws.Range("A62").NumberFormat = "@"
$ws.Range("A62").Value = "2026-01-28"
$ws.Range("A62").Style = "Normal"
$ws.Range("B62").Value = "17:23:43"
$ws.Range("C62").Value = "17:00"
$ws.Range("D62").Value = "Bedroom"
$ws.Range("E62").Value = "In Bed | HR=50 | BR=2"
$ws.Range("F62").Value = "Occupied"

# Row 63
$ws.Range("A63").NumberFormat = "@"
$ws.Range("A63").Value = "2026-01-28"
$ws.Range("A63").Style = "Normal"
$ws.Range("B63").Value = "17:23:44"
$ws.Range("C63").Value = "17:00"
$ws.Range("D63").Value = "Bedroom"
$ws.Range("E63").Value = "In Bed | HR=55 | BR=7"
$ws.Range("F63").Value = "Occupied"

# Row 64
$ws.Range("A64").NumberFormat = "@"
$ws.Range("A64").Value = "2026-01-28"
$ws.Range("A64").Style = "Normal"
$ws.Range("B64").Value = "17:23:45"
$ws.Range("C64").Value = "17:00"
$ws.Range("D64").Value = "Bedroom"
$ws.Range("E64").Value = "In Bed | HR=79 | BR=31"
$ws.Range("F64").Value = "Occupied"

# Row 65
$ws.Range("A65").NumberFormat = "@"
$ws.Range("A65").Value = "2026-01-28"
$ws.Range("A65").Style = "Normal"
$ws.Range("B65").Value = "17:23:46"
$ws.Range("C65").Value = "17:00"
$ws.Range("D65").Value = "Bedroom"
$ws.Range("E65").Value = "In Bed | HR=104 | BR=56"
$ws.Range("F65").Value = "Occupied"

# Row 66
$ws.Range("A66").NumberFormat = "@"
$ws.Range("A66").Value = "2026-01-28"
$ws.Range("A66").Style = "Normal"
$ws.Range("B66").Value = "17:23:47"
$ws.Range("C66").Value = "17:00"
$ws.Range("D66").Value = "Bedroom"
$ws.Range("E66").Value = "In Bed | HR=50 | BR=2"
$ws.Range("F66").Value = "Occupied"

# Row 67
$ws.Range("A67").NumberFormat = "@"
$ws.Range("A67").Value = "2026-01-28"
$ws.Range("A67").Style = "Normal"
$ws.Range("B67").Value = "17:23:49"
$ws.Range("C67").Value = "17:00"
$ws.Range("D67").Value = "Bedroom"
$ws.Range("E67").Value = "In Bed | HR=54 | BR=6"
$ws.Range("F67").Value = "Occupied"

# Row 68
$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "2026-01-28"
$ws.Range("A68").Style = "Normal"
$ws.Range("B68").Value = "17:23:50"
$ws.Range("C68").Value = "17:00"
$ws.Range("D68").Value = "Bedroom"
$ws.Range("E68").Value = "In Bed | HR=50 | BR=2"
$ws.Range("F68").Value = "Occupied"

# Row 69
$ws.Range("A69").NumberFormat = "@"
$ws.Range("A69").Value = "2026-01-28"
$ws.Range("A69").Style = "Normal"
$ws.Range("B69").Value = "17:23:51"
$ws.Range("C69").Value = "17:00"
$ws.Range("D69").Value = "Bedroom"
$ws.Range("E69").Value = "In Bed | HR=116 | BR=68"
$ws.Range("F69").Value = "Occupied"

# Row 70
$ws.Range("A70").NumberFormat = "@"
$ws.Range("A70").Value = "2026-01-28"
$ws.Range("A70").Style = "Normal"
$ws.Range("B70").Value = "17:23:52"
$ws.Range("C70").Value = "17:00"
$ws.Range("D70").Value = "Bedroom"
$ws.Range("E70").Value = "In Bed | HR=92 | BR=44"
$ws.Range("F70").Value = "Occupied"

# Row 71
$ws.Range("A71").NumberFormat = "@"
$ws.Range("A71").Value = "2026-01-28"
$ws.Range("A71").Style = "Normal"
$ws.Range("B71").Value = "17:23:53"
$ws.Range("C71").Value = "17:00"
$ws.Range("D71").Value = "Bedroom"
$ws.Range("E71").Value = "In Bed | HR=102 | BR=54"
$ws.Range("F71").Value = "Occupied"

# Row 72
$ws.Range("A72").NumberFormat = "@"
$ws.Range("A72").Value = "2026-01-28"
$ws.Range("A72").Style = "Normal"
$ws.Range("B72").Value = "17:23:54"
$ws.Range("C72").Value = "17:00"
$ws.Range("D72").Value = "Bedroom"
$ws.Range("E72").Value = "In Bed | HR=81 | BR=33"
$ws.Range("F72").Value = "Occupied"

# Row 73
$ws.Range("A73").NumberFormat = "@"
$ws.Range("A73").Value = "2026-01-28"
$ws.Range("A73").Style = "Normal"
$ws.Range("B73").Value = "17:23:55"
$ws.Range("C73").Value = "17:00"
$ws.Range("D73").Value = "Bedroom"
$ws.Range("E73").Value = "In Bed | HR=50 | BR=2"
$ws.Range("F73").Value = "Occupied"

# Row 74
$ws.Range("A74").NumberFormat = "@"
$ws.Range("A74").Value = "2026-01-28"
$ws.Range("A74").Style = "Normal"
$ws.Range("B74").Value = "17:23:59"
$ws.Range("C74").Value = "17:00"
$ws.Range("D74").Value = "Bedroom"
$ws.Range("E74").Value = "In Bed | HR=111 | BR=63"
$ws.Range("F74").Value = "Occupied"

# Row 75
$ws.Range("A75").NumberFormat = "@"
$ws.Range("A75").Value = "2026-01-28"
$ws.Range("A75").Style = "Normal"
$ws.Range("B75").Value = "17:24:00"
$ws.Range("C75").Value = "17:00"
$ws.Range("D75").Value = "Bedroom"
$ws.Range("E75").Value = "In Bed | HR=99 | BR=51"
$ws.Range("F75").Value = "Occupied"

# Row 76
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "2026-01-28"
$ws.Range("A76").Style = "Normal"
$ws.Range("B76").Value = "17:24:01"
$ws.Range("C76").Value = "17:00"
$ws.Range("D76").Value = "Bedroom"
$ws.Range("E76").Value = "In Bed | HR=96 | BR=48"
$ws.Range("F76").Value = "Occupied"

# Row 77
$ws.Range("A77").NumberFormat = "@"
$ws.Range("A77").Value = "2026-01-28"
$ws.Range("A77").Style = "Normal"
$ws.Range("B77").Value = "17:24:02"
$ws.Range("C77").Value = "17:00"
$ws.Range("D77").Value = "Bedroom"
$ws.Range("E77").Value = "In Bed | HR=108 | BR=60"
$ws.Range("F77").Value = "Occupied"

